# Add a new "2022-Q1" sheet (fund-level holdings) before the "总计" summary
# sheet, and add a corresponding "2022-Q1" row at the top of "总计".

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Create the new "2022-Q1" sheet by cloning "2021-Q4" (same column
#    layout/styles: 基金代码/基金名称/基金规模/股票总仓位/仓位占比/
#    持有市值(亿元)/仓位排名), inserted right before "总计".
# ------------------------------------------------------------------
$template = $wb.Worksheets.Item("2021-Q4")
$zjName = $wb.Worksheets.Item($wb.Worksheets.Count).Name
$zjForPosition = $wb.Worksheets.Item($wb.Worksheets.Count)
$template.Copy($zjForPosition)

$zjAfterCopy = $wb.Worksheets.Item($zjName)
$q1 = $wb.Worksheets.Item($zjAfterCopy.Index - 1)
$q1.Name = "2022-Q1"

# The template ("2021-Q4") has 9 data rows (rows 2-10); the new sheet only
# needs 7 data rows (rows 2-8), so drop the two extra rows.
$q1.Rows.Item(9).Resize(2).Delete()

$fundRows = @(
    @("515210", "国泰中证钢铁ETF", "16.24", "99.25", "3.57", "0.5798", 7),
    @("502023", "鹏华国证钢铁行业指数（LOF）", "15.55", "94.76", "3.49", "0.5427", 7),
    @("168203", "中融国证钢铁行业指数", "4.30", "92.58", "3.39", "0.1458", 7),
    @("013802", "财通资管中证钢铁指数A", "0.11", "90.83", "3.40", "0.0037", 7),
    @("011987", "财通资管智选核心回报6个月持有期混合型发起式证券投资基金A", "0.16", "38.14", "1.36", "0.0022", 5),
    @("013803", "财通资管中证钢铁指数C", "0.02", "90.83", "3.40", "0.0007", 7),
    @("011988", "财通资管智选核心回报6个月持有期混合型发起式证券投资基金C", "0.01", "38.14", "1.36", "0.0001", 5)
)

$r = 2
foreach ($row in $fundRows) {
    $q1.Range("A$r").Value = $r - 2
    $q1.Range("B$r").Value = "'" + $row[0]
    $q1.Range("C$r").Value = $row[1]
    $q1.Range("D$r").Value = "'" + $row[2]
    $q1.Range("E$r").Value = "'" + $row[3]
    $q1.Range("F$r").Value = "'" + $row[4]
    $q1.Range("G$r").Value = "'" + $row[5]
    $q1.Range("H$r").Value = $row[6]
    $r++
}

# ------------------------------------------------------------------
# 2) Insert a new top data row in "总计" for 2022-Q1, shifting the
#    existing history rows down, and renumber the index column.
# ------------------------------------------------------------------
$zj = $wb.Worksheets.Item("总计")
$zj.Rows.Item(2).Insert()
$zj.Range("A2:D2").ClearFormats()

$zj.Range("A2").Value = 0
$zj.Range("B2").Value = "2022-Q1"
$zj.Range("C2").Value = 7
$zj.Range("D2").Value = 1.27

for ($r = 3; $r -le 7; $r++) {
    $zj.Range("A$r").Value = $r - 2
}

# Restore the index-column style (lost on the inserted row) from a
# neighboring row so A2 matches A3:A7.
$zj.Range("A3").Copy()
$zj.Range("A2").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0
